$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dog")
$ws.Activate()

$row = 13

# Copy number formats from row 12 (matches existing style indices) before setting values
$ws.Cells.Item(12, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)
$ws.Cells.Item(12, 3).Copy()
$ws.Cells.Item($row, 3).PasteSpecial(-4122)
$ws.Cells.Item(12, 4).Copy()
$ws.Cells.Item($row, 4).PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item($row, 1).Value = 45803
$ws.Cells.Item($row, 2).Value = "PRESENCE"
$ws.Cells.Item($row, 3).Value = 0.38194444444444442
$ws.Cells.Item($row, 4).Value = 0.4861111111111111
$ws.Cells.Item($row, 5).Value = 15
$ws.Cells.Item($row, 6).Value = 9
$ws.Cells.Item($row, 7).Value = "Overcast, mild"
$ws.Cells.Item($row, 8).Value = $true
$ws.Cells.Item($row, 9).Value = "3 minutes 20 seconds"
$ws.Cells.Item($row, 10).Value = 200
$ws.Cells.Item($row, 11).Value = "Primary sweeps"
$ws.Cells.Item($row, 12).Value = "Worked downhill. Once again on third sweep when Koda picked up odour and belted downhill."

$ws.Range("J14").Select()
# Scroll so column E becomes the leftmost visible column (mirrors topLeftCell="E1")
$excel.ActiveWindow.ScrollColumn = 5
